$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2; existing rows 2-22 shift down to 3-23 with all
# their values/styles intact.
$ws.Rows("2:2").Insert()

# The row that used to be row 2 is now row 3 - bump a few of its values.
$ws.Range("C3").Value = 418324115
$ws.Range("E3").Value = 18329953645
$ws.Range("F3").Value = "https://raw.githubusercontent.com/Aishee002/ADUserdata/main/ADUserdata.xlsx"
$ws.Range("G3").Value = "1/26/2026"

# Populate the new row 2 with the promo / QQ-group entry.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "交流Q群"
$ws.Range("C2").Value = 321672464
$ws.Range("F2").Value = "https://qm.qq.com/q/rSlKgUAWZO"
$ws.Range("G2").Value = "9/9/2099"
$ws.Range("G2").NumberFormat = "m/d/yy"

# Style the new row with a green fill (matches the workbook's accent-6 theme color).
$ws.Range("A2:F2").Interior.Color = 4697456
$ws.Range("G2").Interior.Color = 4697456

# Selection moved in the saved file.
$ws.Range("I5").Select()
